# edit.ps1
# Applies the "typos im Antrag, fixing #99" commit to the Tempo30 application template.
$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Resize the two tiny decorative "fold mark" drawings (docPr id 1 & 2)
#    that live as inline shapes inside the very first paragraph of the
#    document. Word's InlineShape.Width/Height only rewrite the outer
#    <wp:extent>, while the diff also needs the inner <a:xfrm><a:ext>
#    of the custom-geometry shape updated to a *different* value, so we
#    replace the whole first paragraph (which also carries the "Receiver
#    Address" anchored text box) with a copy that has both pairs of
#    numbers patched.
# -----------------------------------------------------------------
$para1Xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:pStyle w:val="Date"/></w:pPr><w:r><w:rPr/><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0"><wp:extent cx="274955" cy="5080"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:docPr id="1" name=""/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="274320" cy="4320"/></a:xfrm><a:custGeom><a:avLst/><a:gdLst/><a:ahLst/><a:rect l="0" t="0" r="r" b="b"/><a:pathLst><a:path w="0" h="0"/></a:pathLst></a:custGeom><a:noFill/><a:ln><a:noFill/></a:ln></wps:spPr><wps:style><a:lnRef idx="0"/><a:fillRef idx="0"/><a:effectRef idx="0"/><a:fontRef idx="minor"/></wps:style><wps:bodyPr/></wps:wsp></a:graphicData></a:graphic></wp:inline></w:drawing></mc:Choice><mc:Fallback><w:pict/></mc:Fallback></mc:AlternateContent></w:r><w:r><w:rPr/><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0"><wp:extent cx="274955" cy="5080"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:docPr id="2" name=""/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="274320" cy="4320"/></a:xfrm><a:custGeom><a:avLst/><a:gdLst/><a:ahLst/><a:rect l="0" t="0" r="r" b="b"/><a:pathLst><a:path w="0" h="0"/></a:pathLst></a:custGeom><a:noFill/><a:ln><a:noFill/></a:ln></wps:spPr><wps:style><a:lnRef idx="0"/><a:fillRef idx="0"/><a:effectRef idx="0"/><a:fontRef idx="minor"/></wps:style><wps:bodyPr/></wps:wsp></a:graphicData></a:graphic></wp:inline></w:drawing></mc:Choice><mc:Fallback><w:pict/></mc:Fallback></mc:AlternateContent></w:r><w:r><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor behindDoc="0" distT="0" distB="0" distL="0" distR="0" simplePos="0" locked="0" layoutInCell="1" allowOverlap="1" relativeHeight="2"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="page"><wp:posOffset>760095</wp:posOffset></wp:positionH><wp:positionV relativeFrom="page"><wp:posOffset>1805305</wp:posOffset></wp:positionV><wp:extent cx="2680335" cy="919480"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapNone/><wp:docPr id="3" name="Receiver Address"/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr txBox="1"/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="2680335" cy="919480"/></a:xfrm><a:prstGeom prst="rect"/></wps:spPr><wps:txbx><w:txbxContent><w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr><w:sz w:val="22"/><w:spacing w:val="2"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rFonts w:ascii="Cumberland" w:hAnsi="Cumberland" w:eastAsia="Andale Sans UI" w:cs="Tahoma"/><w:color w:val="00000A"/><w:lang w:val="de-DE" w:eastAsia="de-DE" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr/><w:t>An das</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr><w:sz w:val="22"/><w:spacing w:val="2"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rFonts w:ascii="Cumberland" w:hAnsi="Cumberland" w:eastAsia="Andale Sans UI" w:cs="Tahoma"/><w:color w:val="00000A"/><w:lang w:val="de-DE" w:eastAsia="de-DE" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr/><w:t>{Polizei}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr><w:sz w:val="22"/><w:spacing w:val="2"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rFonts w:ascii="Cumberland" w:hAnsi="Cumberland" w:eastAsia="Andale Sans UI" w:cs="Tahoma"/><w:color w:val="00000A"/><w:lang w:val="de-DE" w:eastAsia="de-DE" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr/><w:t>- Straßenverkehrsbehörde -</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr><w:sz w:val="22"/><w:spacing w:val="2"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rFonts w:ascii="Cumberland" w:hAnsi="Cumberland" w:eastAsia="Andale Sans UI" w:cs="Tahoma"/><w:color w:val="00000A"/><w:lang w:val="de-DE" w:eastAsia="de-DE" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr/><w:t>{PolizeiStr}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:before="0" w:after="120"/><w:rPr><w:sz w:val="22"/><w:spacing w:val="2"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rFonts w:ascii="Cumberland" w:hAnsi="Cumberland" w:eastAsia="Andale Sans UI" w:cs="Tahoma"/><w:color w:val="00000A"/><w:lang w:val="de-DE" w:eastAsia="de-DE" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr/><w:t>{PolizeiPLZ} {PolizeiOrt}</w:t></w:r></w:p></w:txbxContent></wps:txbx><wps:bodyPr anchor="t" lIns="0" tIns="0" rIns="0" bIns="0"><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:rect stroked="f" strokeweight="0pt" style="position:absolute;width:211.05pt;height:72.4pt;mso-wrap-distance-left:0pt;mso-wrap-distance-right:0pt;mso-wrap-distance-top:0pt;mso-wrap-distance-bottom:0pt;margin-top:142.15pt;mso-position-vertical-relative:page;margin-left:59.85pt;mso-position-horizontal-relative:page"><v:textbox inset="0in,0in,0in,0in"><w:txbxContent><w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr><w:sz w:val="22"/><w:spacing w:val="2"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rFonts w:ascii="Cumberland" w:hAnsi="Cumberland" w:eastAsia="Andale Sans UI" w:cs="Tahoma"/><w:color w:val="00000A"/><w:lang w:val="de-DE" w:eastAsia="de-DE" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr/><w:t>An das</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr><w:sz w:val="22"/><w:spacing w:val="2"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rFonts w:ascii="Cumberland" w:hAnsi="Cumberland" w:eastAsia="Andale Sans UI" w:cs="Tahoma"/><w:color w:val="00000A"/><w:lang w:val="de-DE" w:eastAsia="de-DE" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr/><w:t>{Polizei}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr><w:sz w:val="22"/><w:spacing w:val="2"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rFonts w:ascii="Cumberland" w:hAnsi="Cumberland" w:eastAsia="Andale Sans UI" w:cs="Tahoma"/><w:color w:val="00000A"/><w:lang w:val="de-DE" w:eastAsia="de-DE" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr/><w:t>- Straßenverkehrsbehörde -</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr><w:sz w:val="22"/><w:spacing w:val="2"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rFonts w:ascii="Cumberland" w:hAnsi="Cumberland" w:eastAsia="Andale Sans UI" w:cs="Tahoma"/><w:color w:val="00000A"/><w:lang w:val="de-DE" w:eastAsia="de-DE" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr/><w:t>{PolizeiStr}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:before="0" w:after="120"/><w:rPr><w:sz w:val="22"/><w:spacing w:val="2"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rFonts w:ascii="Cumberland" w:hAnsi="Cumberland" w:eastAsia="Andale Sans UI" w:cs="Tahoma"/><w:color w:val="00000A"/><w:lang w:val="de-DE" w:eastAsia="de-DE" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr/><w:t>{PolizeiPLZ} {PolizeiOrt}</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect></w:pict></mc:Fallback></mc:AlternateContent></w:r><w:r><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor behindDoc="0" distT="0" distB="0" distL="0" distR="0" simplePos="0" locked="0" layoutInCell="1" allowOverlap="1" relativeHeight="3"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>50165</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>41910</wp:posOffset></wp:positionV><wp:extent cx="2670175" cy="188595"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapNone/><wp:docPr id="4" name="Sender Address Repeated"/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr txBox="1"/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="2670175" cy="188595"/></a:xfrm><a:prstGeom prst="rect"/></wps:spPr><wps:txbx><w:txbxContent><w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:i w:val="false"/><w:iCs w:val="false"/><w:sz w:val="13"/><w:szCs w:val="13"/></w:rPr><w:t>{Name}, {AdrStr}, {AdrPLZ} {AdrOrt}</w:t></w:r></w:p></w:txbxContent></wps:txbx><wps:bodyPr anchor="t" lIns="0" tIns="0" rIns="0" bIns="0"><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:rect stroked="f" strokeweight="0pt" style="position:absolute;width:210.25pt;height:14.85pt;mso-wrap-distance-left:0pt;mso-wrap-distance-right:0pt;mso-wrap-distance-top:0pt;mso-wrap-distance-bottom:0pt;margin-top:3.3pt;mso-position-vertical-relative:text;margin-left:3.95pt;mso-position-horizontal-relative:text"><v:textbox inset="0in,0in,0in,0in"><w:txbxContent><w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:i w:val="false"/><w:iCs w:val="false"/><w:sz w:val="13"/><w:szCs w:val="13"/></w:rPr><w:t>{Name}, {AdrStr}, {AdrPLZ} {AdrOrt}</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect></w:pict></mc:Fallback></mc:AlternateContent></w:r><w:r><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor behindDoc="0" distT="0" distB="0" distL="0" distR="0" simplePos="0" locked="0" layoutInCell="1" allowOverlap="1" relativeHeight="4"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>3203575</wp:posOffset></wp:positionH><wp:positionV relativeFrom="page"><wp:posOffset>590550</wp:posOffset></wp:positionV><wp:extent cx="2953385" cy="829310"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapNone/><wp:docPr id="5" name="Sender Address"/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr txBox="1"/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="2953385" cy="829310"/></a:xfrm><a:prstGeom prst="rect"/></wps:spPr><wps:txbx><w:txbxContent><w:p><w:pPr><w:pStyle w:val="Berschrift"/><w:spacing w:before="0" w:after="120"/><w:rPr><w:sz w:val="28"/><w:spacing w:val="2"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rFonts w:ascii="Courier" w:hAnsi="Courier" w:eastAsia="Droid Sans Fallback" w:cs="FreeSans"/><w:color w:val="00000A"/><w:lang w:val="de-DE" w:eastAsia="de-DE" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr/><w:t>{Name}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Kopfzeilerechts"/><w:rPr><w:sz w:val="18"/><w:spacing w:val="2"/><w:sz w:val="18"/><w:szCs w:val="24"/><w:rFonts w:ascii="Courier" w:hAnsi="Courier" w:eastAsia="Andale Sans UI" w:cs="Tahoma"/><w:color w:val="00000A"/><w:lang w:val="de-DE" w:eastAsia="de-DE" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr/><w:t>{AdrStr}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Kopfzeilerechts"/><w:shd w:val="clear" w:color="000000" w:themeColor="" w:themeTint="0" w:themeShade="0" w:fill="FFFFFF" w:themeFill="" w:themeFillTint="0" w:themeFillShade="0"/><w:rPr><w:sz w:val="18"/><w:spacing w:val="2"/><w:sz w:val="18"/><w:szCs w:val="24"/><w:rFonts w:ascii="Courier" w:hAnsi="Courier" w:eastAsia="Andale Sans UI" w:cs="Tahoma"/><w:color w:val="00000A"/><w:lang w:val="de-DE" w:eastAsia="de-DE" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr/><w:t>{AdrPLZ} {AdrOrt}</w:t></w:r></w:p></w:txbxContent></wps:txbx><wps:bodyPr anchor="t" lIns="0" tIns="0" rIns="0" bIns="0"><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:rect stroked="f" strokeweight="0pt" style="position:absolute;width:232.55pt;height:65.3pt;mso-wrap-distance-left:0pt;mso-wrap-distance-right:0pt;mso-wrap-distance-top:0pt;mso-wrap-distance-bottom:0pt;margin-top:46.5pt;mso-position-vertical-relative:page;margin-left:252.25pt;mso-position-horizontal-relative:text"><v:textbox inset="0in,0in,0in,0in"><w:txbxContent><w:p><w:pPr><w:pStyle w:val="Berschrift"/><w:spacing w:before="0" w:after="120"/><w:rPr><w:sz w:val="28"/><w:spacing w:val="2"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:rFonts w:ascii="Courier" w:hAnsi="Courier" w:eastAsia="Droid Sans Fallback" w:cs="FreeSans"/><w:color w:val="00000A"/><w:lang w:val="de-DE" w:eastAsia="de-DE" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr/><w:t>{Name}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Kopfzeilerechts"/><w:rPr><w:sz w:val="18"/><w:spacing w:val="2"/><w:sz w:val="18"/><w:szCs w:val="24"/><w:rFonts w:ascii="Courier" w:hAnsi="Courier" w:eastAsia="Andale Sans UI" w:cs="Tahoma"/><w:color w:val="00000A"/><w:lang w:val="de-DE" w:eastAsia="de-DE" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr/><w:t>{AdrStr}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Kopfzeilerechts"/><w:shd w:val="clear" w:color="000000" w:themeColor="" w:themeTint="0" w:themeShade="0" w:fill="FFFFFF" w:themeFill="" w:themeFillTint="0" w:themeFillShade="0"/><w:rPr><w:sz w:val="18"/><w:spacing w:val="2"/><w:sz w:val="18"/><w:szCs w:val="24"/><w:rFonts w:ascii="Courier" w:hAnsi="Courier" w:eastAsia="Andale Sans UI" w:cs="Tahoma"/><w:color w:val="00000A"/><w:lang w:val="de-DE" w:eastAsia="de-DE" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr/><w:t>{AdrPLZ} {AdrOrt}</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect></w:pict></mc:Fallback></mc:AlternateContent></w:r><w:r><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor behindDoc="0" distT="0" distB="0" distL="0" distR="0" simplePos="0" locked="0" layoutInCell="1" allowOverlap="1" relativeHeight="5"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="page"><wp:posOffset>241300</wp:posOffset></wp:positionH><wp:positionV relativeFrom="page"><wp:posOffset>3420110</wp:posOffset></wp:positionV><wp:extent cx="360045" cy="266065"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapNone/><wp:docPr id="6" name="Bend Marks"/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr txBox="1"/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="360045" cy="266065"/></a:xfrm><a:prstGeom prst="rect"/></wps:spPr><wps:txbx><w:txbxContent><w:p><w:pPr><w:pStyle w:val="Rahmeninhalt"/><w:spacing w:before="0" w:after="0"/><w:rPr><w:sz w:val="22"/><w:spacing w:val="2"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rFonts w:ascii="Cumberland" w:hAnsi="Cumberland" w:eastAsia="Andale Sans UI" w:cs="Tahoma"/><w:color w:val="00000A"/><w:lang w:val="de-DE" w:eastAsia="de-DE" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr/></w:r></w:p></w:txbxContent></wps:txbx><wps:bodyPr anchor="t" lIns="53975" tIns="53975" rIns="53975" bIns="53975"><a:spAutoFit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:rect stroked="f" strokeweight="0pt" style="position:absolute;width:28.35pt;height:20.95pt;mso-wrap-distance-left:0pt;mso-wrap-distance-right:0pt;mso-wrap-distance-top:0pt;mso-wrap-distance-bottom:0pt;margin-top:269.3pt;mso-position-vertical-relative:page;margin-left:19pt;mso-position-horizontal-relative:page"><v:textbox inset="0.0590277777777778in,0.0590277777777778in,0.0590277777777778in,0.0590277777777778in"><w:txbxContent><w:p><w:pPr><w:pStyle w:val="Rahmeninhalt"/><w:spacing w:before="0" w:after="0"/><w:rPr><w:sz w:val="22"/><w:spacing w:val="2"/><w:sz w:val="22"/><w:szCs w:val="24"/><w:rFonts w:ascii="Cumberland" w:hAnsi="Cumberland" w:eastAsia="Andale Sans UI" w:cs="Tahoma"/><w:color w:val="00000A"/><w:lang w:val="de-DE" w:eastAsia="de-DE" w:bidi="ar-SA"/></w:rPr></w:pPr><w:r><w:rPr/></w:r></w:p></w:txbxContent></v:textbox></v:rect></w:pict></mc:Fallback></mc:AlternateContent></w:r><w:r/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.InsertXML($para1Xml)

# -----------------------------------------------------------------
# 2) Simple literal text fixes (Find/Replace, one unique match each)
# -----------------------------------------------------------------
$d.Content.Find.Execute("21.11.16", $true, $false, $false, $false, $false, $true, 1, $false, "23.11.16", 2) | Out-Null

$d.Content.Find.Execute("[ ] für Kinder werden aufgrund der hohen Geschwindigkeit gefährlich", $true, $false, $false, $false, $false, $true, 1, $false, "[ ] für Kinder aufgrund der hohen Geschwindigkeit gefährlich", 2) | Out-Null

$d.Content.Find.Execute("[ ] durch  häufiges Parken in zweiter Reihe geprägt", $true, $false, $false, $false, $false, $true, 1, $false, "[ ] durch häufiges Parken in zweiter Reihe geprägt", 2) | Out-Null

$d.Content.Find.Execute("[ ] Häufiges Anfahren und Abbremsen geprägt (z.B. Parkplatzsuche etc.)", $true, $false, $false, $false, $false, $true, 1, $false, "[ ] durch häufiges Anfahren und Abbremsen geprägt (z.B. Parkplatzsuche etc.)", 2) | Out-Null

$d.Content.Find.Execute("unserem Hausam Tag", $true, $false, $false, $false, $false, $true, 1, $false, "unserem Haus am Tag", 2) | Out-Null

$d.Content.Find.Execute("anpassen. Desto mehr Details", $true, $false, $false, $false, $false, $true, 1, $false, "anpassen. Je mehr Details", 2) | Out-Null

$d.Content.Find.Execute("wenn Sie dem ADFC-Hamburg informieren", $true, $false, $false, $false, $false, $true, 1, $false, "wenn Sie den ADFC-Hamburg informieren", 2) | Out-Null

$d.Content.Find.Execute("Halten Sie dem ADFC auf dem Laufenden", $true, $false, $false, $false, $false, $true, 1, $false, "Halten Sie den ADFC auf dem Laufenden", 2) | Out-Null

$d.Content.Find.Execute("Bei abgelehnten Antrag oder keiner Antwort", $true, $false, $false, $false, $false, $true, 1, $false, "Bei abgelehntem Antrag oder keiner Antwort", 2) | Out-Null

# -----------------------------------------------------------------
# 3) Insert the new checklist paragraph right after the
#    "... Anfahren und Abbremsen ..." bullet.
# -----------------------------------------------------------------
$target = "[ ] durch häufiges Anfahren und Abbremsen geprägt (z.B. Parkplatzsuche etc.)"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $target) {
        $p.Range.InsertParagraphAfter()
        $newPara = $d.Paragraphs.Item($i + 1)
        $newPara.Range.Text = "[ ] durch zahlreiche Konflikte zwischen Radfahrern, Autofahrern und Fußgängern geprägt"
        break
    }
}
